# "Add capacity construction for grid reliability in every hour" (#232)
#
# On the "About" sheet, a new explanatory note is inserted above the
# existing "shareweights" notes block:  two rows are inserted at row 10
# (pushing the rest of the notes down by two rows - one row for the new
# text, one blank spacer row, matching the existing blank-row spacing
# pattern used throughout this sheet), and the new row 10 is filled with
# bold text, styled like the other section headers on this sheet (e.g.
# A1, A5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert two blank rows at row 10; everything from the old row 10 down
# shifts to row 12+.
$ws.Range("A10:A11").EntireRow.Insert()

# Populate the newly-inserted row 10 with the new note, bolded.
$ws.Range("A10").Value = "The EPS assumes shareweights will be between 0 and 1 (inclusive)."
$ws.Range("A10").Font.Bold = $true

# The sheet is now viewed a bit more zoomed-in.
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
